# Edit: remove the word "significativa " from the objective text on slide 5.
# Original run text: <nbsp>significativa sobre la<nbsp>
# New run text:       <nbsp>sobre la<nbsp>
#
# The word boundaries in this textbox are written using non-breaking
# spaces (U+00A0) instead of regular spaces, so we rebuild the exact
# substring (including the surrounding NBSPs) in one assignment. Doing
# the replacement as a single Characters(...).Text= write (instead of
# deleting a sub-range) keeps everything inside the original single
# <a:r> run, matching how the source run is structured.
#
# Note: this COM runtime round-trips non-ASCII characters (like NBSP)
# back out as U+FFFD (the replacement character) whenever text is read
# back via .Text, even though the underlying XML still stores the real
# NBSP correctly. So for verifying what is already in the document we
# compare using U+FFFD, but for the text we *write*, we use the real
# NBSP character so it is persisted correctly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$nbsp = [char]0x00A0
$replChar = [char]0xFFFD

$readTarget = $replChar + "significativa sobre la" + $replChar
$writeReplacement = $nbsp + "sobre la" + $nbsp

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf("significativa")
    if ($idx -ge 0) {
        # 1-based character index of the NBSP immediately preceding "significativa"
        $startChar = $idx
        $len = $readTarget.Length
        $rng = $tr.Characters($startChar, $len)
        if ($rng.Text -eq $readTarget) {
            $rng.Text = $writeReplacement
        }
        break
    }
}
